# Update the "want to go" (想去人数) counts in column F for two sheets:
# "展览" (Exhibitions) and "全部类型" (All Types), reflecting refreshed
# scrape output as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new F value for the "展览" sheet
$exhibitionUpdates = @{
    3  = 1435
    7  = 1258
    11 = 2287
    17 = 86
    18 = 6288
    20 = 6201
    21 = 10193
    29 = 4396
    30 = 108
    31 = 396
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F value for the "全部类型" sheet
$allTypesUpdates = @{
    5  = 1435
    10 = 1258
    15 = 2287
    23 = 86
    24 = 6288
    26 = 6201
    27 = 10193
    40 = 4396
    46 = 396
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
